# Auto-generated Word COM-interop script
# Applies sequential text replacements (date + 100 arithmetic problems)
# Order is topologically sorted so that a replacement's output never
# collides with a not-yet-processed replacement's search text.

$d = $word.ActiveDocument

$count = 0
if ($d.Content.Find.Execute("2024-01-04 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-01-05 Friday", 2)) { $count++ }
if ($d.Content.Find.Execute("95-89=6", $true, $false, $false, $false, $false, $true, 1, $false, "19+35=54", 2)) { $count++ }
if ($d.Content.Find.Execute("28+57=85", $true, $false, $false, $false, $false, $true, 1, $false, "35-27=8", 2)) { $count++ }
if ($d.Content.Find.Execute("43-37=6", $true, $false, $false, $false, $false, $true, 1, $false, "71-15=56", 2)) { $count++ }
if ($d.Content.Find.Execute("92-56=36", $true, $false, $false, $false, $false, $true, 1, $false, "91-89=2", 2)) { $count++ }
if ($d.Content.Find.Execute("82-34=48", $true, $false, $false, $false, $false, $true, 1, $false, "38+46=84", 2)) { $count++ }
if ($d.Content.Find.Execute("84-38=46", $true, $false, $false, $false, $false, $true, 1, $false, "58+36=94", 2)) { $count++ }
if ($d.Content.Find.Execute("80-69=11", $true, $false, $false, $false, $false, $true, 1, $false, "55+8=63", 2)) { $count++ }
if ($d.Content.Find.Execute("49+48=97", $true, $false, $false, $false, $false, $true, 1, $false, "86+8=94", 2)) { $count++ }
if ($d.Content.Find.Execute("6+5=11", $true, $false, $false, $false, $false, $true, 1, $false, "44+48=92", 2)) { $count++ }
if ($d.Content.Find.Execute("47+38=85", $true, $false, $false, $false, $false, $true, 1, $false, "9+75=84", 2)) { $count++ }
if ($d.Content.Find.Execute("93-79=14", $true, $false, $false, $false, $false, $true, 1, $false, "25+49=74", 2)) { $count++ }
if ($d.Content.Find.Execute("70-48=22", $true, $false, $false, $false, $false, $true, 1, $false, "16+67=83", 2)) { $count++ }
if ($d.Content.Find.Execute("9+33=42", $true, $false, $false, $false, $false, $true, 1, $false, "37+8=45", 2)) { $count++ }
if ($d.Content.Find.Execute("28-19=9", $true, $false, $false, $false, $false, $true, 1, $false, "17+44=61", 2)) { $count++ }
if ($d.Content.Find.Execute("41-22=19", $true, $false, $false, $false, $false, $true, 1, $false, "39+55=94", 2)) { $count++ }
if ($d.Content.Find.Execute("13+49=62", $true, $false, $false, $false, $false, $true, 1, $false, "74-47=27", 2)) { $count++ }
if ($d.Content.Find.Execute("85-39=46", $true, $false, $false, $false, $false, $true, 1, $false, "3+9=12", 2)) { $count++ }
if ($d.Content.Find.Execute("52+39=91", $true, $false, $false, $false, $false, $true, 1, $false, "36+15=51", 2)) { $count++ }
if ($d.Content.Find.Execute("79+17=96", $true, $false, $false, $false, $false, $true, 1, $false, "40-3=37", 2)) { $count++ }
if ($d.Content.Find.Execute("82-15=67", $true, $false, $false, $false, $false, $true, 1, $false, "38+47=85", 2)) { $count++ }
if ($d.Content.Find.Execute("73-36=37", $true, $false, $false, $false, $false, $true, 1, $false, "4+59=63", 2)) { $count++ }
if ($d.Content.Find.Execute("55+27=82", $true, $false, $false, $false, $false, $true, 1, $false, "15-7=8", 2)) { $count++ }
if ($d.Content.Find.Execute("83-48=35", $true, $false, $false, $false, $false, $true, 1, $false, "7+69=76", 2)) { $count++ }
if ($d.Content.Find.Execute("46-37=9", $true, $false, $false, $false, $false, $true, 1, $false, "83-48=35", 2)) { $count++ }
if ($d.Content.Find.Execute("80-57=23", $true, $false, $false, $false, $false, $true, 1, $false, "92-54=38", 2)) { $count++ }
if ($d.Content.Find.Execute("73-27=46", $true, $false, $false, $false, $false, $true, 1, $false, "20-6=14", 2)) { $count++ }
if ($d.Content.Find.Execute("80-48=32", $true, $false, $false, $false, $false, $true, 1, $false, "9+49=58", 2)) { $count++ }
if ($d.Content.Find.Execute("37+36=73", $true, $false, $false, $false, $false, $true, 1, $false, "21-13=8", 2)) { $count++ }
if ($d.Content.Find.Execute("19+75=94", $true, $false, $false, $false, $false, $true, 1, $false, "92-15=77", 2)) { $count++ }
if ($d.Content.Find.Execute("60-56=4", $true, $false, $false, $false, $false, $true, 1, $false, "28+36=64", 2)) { $count++ }
if ($d.Content.Find.Execute("82-26=56", $true, $false, $false, $false, $false, $true, 1, $false, "70-4=66", 2)) { $count++ }
if ($d.Content.Find.Execute("68-39=29", $true, $false, $false, $false, $false, $true, 1, $false, "59+8=67", 2)) { $count++ }
if ($d.Content.Find.Execute("41-35=6", $true, $false, $false, $false, $false, $true, 1, $false, "54-9=45", 2)) { $count++ }
if ($d.Content.Find.Execute("17+48=65", $true, $false, $false, $false, $false, $true, 1, $false, "16+8=24", 2)) { $count++ }
if ($d.Content.Find.Execute("36-29=7", $true, $false, $false, $false, $false, $true, 1, $false, "65+18=83", 2)) { $count++ }
if ($d.Content.Find.Execute("72-8=64", $true, $false, $false, $false, $false, $true, 1, $false, "91-13=78", 2)) { $count++ }
if ($d.Content.Find.Execute("33-6=27", $true, $false, $false, $false, $false, $true, 1, $false, "79+14=93", 2)) { $count++ }
if ($d.Content.Find.Execute("59+33=92", $true, $false, $false, $false, $false, $true, 1, $false, "56-49=7", 2)) { $count++ }
if ($d.Content.Find.Execute("94-59=35", $true, $false, $false, $false, $false, $true, 1, $false, "81-29=52", 2)) { $count++ }
if ($d.Content.Find.Execute("38+53=91", $true, $false, $false, $false, $false, $true, 1, $false, "5+39=44", 2)) { $count++ }
if ($d.Content.Find.Execute("83-64=19", $true, $false, $false, $false, $false, $true, 1, $false, "74+19=93", 2)) { $count++ }
if ($d.Content.Find.Execute("79+13=92", $true, $false, $false, $false, $false, $true, 1, $false, "59+19=78", 2)) { $count++ }
if ($d.Content.Find.Execute("80-51=29", $true, $false, $false, $false, $false, $true, 1, $false, "82-7=75", 2)) { $count++ }
if ($d.Content.Find.Execute("12-5=7", $true, $false, $false, $false, $false, $true, 1, $false, "38-29=9", 2)) { $count++ }
if ($d.Content.Find.Execute("81-42=39", $true, $false, $false, $false, $false, $true, 1, $false, "56-39=17", 2)) { $count++ }
if ($d.Content.Find.Execute("72-24=48", $true, $false, $false, $false, $false, $true, 1, $false, "3+48=51", 2)) { $count++ }
if ($d.Content.Find.Execute("10-2=8", $true, $false, $false, $false, $false, $true, 1, $false, "61-2=59", 2)) { $count++ }
if ($d.Content.Find.Execute("36+49=85", $true, $false, $false, $false, $false, $true, 1, $false, "64-45=19", 2)) { $count++ }
if ($d.Content.Find.Execute("26-18=8", $true, $false, $false, $false, $false, $true, 1, $false, "63-27=36", 2)) { $count++ }
if ($d.Content.Find.Execute("85-8=77", $true, $false, $false, $false, $false, $true, 1, $false, "44+49=93", 2)) { $count++ }
if ($d.Content.Find.Execute("18+55=73", $true, $false, $false, $false, $false, $true, 1, $false, "36+16=52", 2)) { $count++ }
if ($d.Content.Find.Execute("62-55=7", $true, $false, $false, $false, $false, $true, 1, $false, "55-17=38", 2)) { $count++ }
if ($d.Content.Find.Execute("93-57=36", $true, $false, $false, $false, $false, $true, 1, $false, "27+66=93", 2)) { $count++ }
if ($d.Content.Find.Execute("36+8=44", $true, $false, $false, $false, $false, $true, 1, $false, "57+9=66", 2)) { $count++ }
if ($d.Content.Find.Execute("90-9=81", $true, $false, $false, $false, $false, $true, 1, $false, "71-54=17", 2)) { $count++ }
if ($d.Content.Find.Execute("87-39=48", $true, $false, $false, $false, $false, $true, 1, $false, "65-16=49", 2)) { $count++ }
if ($d.Content.Find.Execute("48+48=96", $true, $false, $false, $false, $false, $true, 1, $false, "5+66=71", 2)) { $count++ }
if ($d.Content.Find.Execute("41-24=17", $true, $false, $false, $false, $false, $true, 1, $false, "90-24=66", 2)) { $count++ }
if ($d.Content.Find.Execute("58+23=81", $true, $false, $false, $false, $false, $true, 1, $false, "97-9=88", 2)) { $count++ }
if ($d.Content.Find.Execute("92-35=57", $true, $false, $false, $false, $false, $true, 1, $false, "6+6=12", 2)) { $count++ }
if ($d.Content.Find.Execute("25-19=6", $true, $false, $false, $false, $false, $true, 1, $false, "82-13=69", 2)) { $count++ }
if ($d.Content.Find.Execute("44-37=7", $true, $false, $false, $false, $false, $true, 1, $false, "59+9=68", 2)) { $count++ }
if ($d.Content.Find.Execute("14+57=71", $true, $false, $false, $false, $false, $true, 1, $false, "85-69=16", 2)) { $count++ }
if ($d.Content.Find.Execute("36+58=94", $true, $false, $false, $false, $false, $true, 1, $false, "6+39=45", 2)) { $count++ }
if ($d.Content.Find.Execute("72-9=63", $true, $false, $false, $false, $false, $true, 1, $false, "96-17=79", 2)) { $count++ }
if ($d.Content.Find.Execute("7+26=33", $true, $false, $false, $false, $false, $true, 1, $false, "94-87=7", 2)) { $count++ }
if ($d.Content.Find.Execute("32+9=41", $true, $false, $false, $false, $false, $true, 1, $false, "64-5=59", 2)) { $count++ }
if ($d.Content.Find.Execute("87+9=96", $true, $false, $false, $false, $false, $true, 1, $false, "8+17=25", 2)) { $count++ }
if ($d.Content.Find.Execute("18+29=47", $true, $false, $false, $false, $false, $true, 1, $false, "49+49=98", 2)) { $count++ }
if ($d.Content.Find.Execute("84-25=59", $true, $false, $false, $false, $false, $true, 1, $false, "92-53=39", 2)) { $count++ }
if ($d.Content.Find.Execute("91-65=26", $true, $false, $false, $false, $false, $true, 1, $false, "73+8=81", 2)) { $count++ }
if ($d.Content.Find.Execute("9+46=55", $true, $false, $false, $false, $false, $true, 1, $false, "28+18=46", 2)) { $count++ }
if ($d.Content.Find.Execute("17+4=21", $true, $false, $false, $false, $false, $true, 1, $false, "37+56=93", 2)) { $count++ }
if ($d.Content.Find.Execute("48+24=72", $true, $false, $false, $false, $false, $true, 1, $false, "83-9=74", 2)) { $count++ }
if ($d.Content.Find.Execute("92-66=26", $true, $false, $false, $false, $false, $true, 1, $false, "47+16=63", 2)) { $count++ }
if ($d.Content.Find.Execute("62-59=3", $true, $false, $false, $false, $false, $true, 1, $false, "80-27=53", 2)) { $count++ }
if ($d.Content.Find.Execute("49+22=71", $true, $false, $false, $false, $false, $true, 1, $false, "82-47=35", 2)) { $count++ }
if ($d.Content.Find.Execute("88-69=19", $true, $false, $false, $false, $false, $true, 1, $false, "24-9=15", 2)) { $count++ }
if ($d.Content.Find.Execute("23+18=41", $true, $false, $false, $false, $false, $true, 1, $false, "55+37=92", 2)) { $count++ }
if ($d.Content.Find.Execute("81-12=69", $true, $false, $false, $false, $false, $true, 1, $false, "53-49=4", 2)) { $count++ }
if ($d.Content.Find.Execute("17+45=62", $true, $false, $false, $false, $false, $true, 1, $false, "71-26=45", 2)) { $count++ }
if ($d.Content.Find.Execute("9+58=67", $true, $false, $false, $false, $false, $true, 1, $false, "48+38=86", 2)) { $count++ }
if ($d.Content.Find.Execute("81-47=34", $true, $false, $false, $false, $false, $true, 1, $false, "58+3=61", 2)) { $count++ }
if ($d.Content.Find.Execute("61-49=12", $true, $false, $false, $false, $false, $true, 1, $false, "62-58=4", 2)) { $count++ }
if ($d.Content.Find.Execute("53+9=62", $true, $false, $false, $false, $false, $true, 1, $false, "52-18=34", 2)) { $count++ }
if ($d.Content.Find.Execute("15+29=44", $true, $false, $false, $false, $false, $true, 1, $false, "90-36=54", 2)) { $count++ }
if ($d.Content.Find.Execute("71-57=14", $true, $false, $false, $false, $false, $true, 1, $false, "60-21=39", 2)) { $count++ }
if ($d.Content.Find.Execute("65-6=59", $true, $false, $false, $false, $false, $true, 1, $false, "91-17=74", 2)) { $count++ }
if ($d.Content.Find.Execute("77-58=19", $true, $false, $false, $false, $false, $true, 1, $false, "85-26=59", 2)) { $count++ }
if ($d.Content.Find.Execute("25-17=8", $true, $false, $false, $false, $false, $true, 1, $false, "72-3=69", 2)) { $count++ }
if ($d.Content.Find.Execute("74-67=7", $true, $false, $false, $false, $false, $true, 1, $false, "9+3=12", 2)) { $count++ }
if ($d.Content.Find.Execute("7+4=11", $true, $false, $false, $false, $false, $true, 1, $false, "85-47=38", 2)) { $count++ }
if ($d.Content.Find.Execute("67-58=9", $true, $false, $false, $false, $false, $true, 1, $false, "75+8=83", 2)) { $count++ }
if ($d.Content.Find.Execute("21-5=16", $true, $false, $false, $false, $false, $true, 1, $false, "25-18=7", 2)) { $count++ }
if ($d.Content.Find.Execute("32-7=25", $true, $false, $false, $false, $false, $true, 1, $false, "56+9=65", 2)) { $count++ }
if ($d.Content.Find.Execute("80-59=21", $true, $false, $false, $false, $false, $true, 1, $false, "4+59=63", 2)) { $count++ }
if ($d.Content.Find.Execute("37-18=19", $true, $false, $false, $false, $false, $true, 1, $false, "51-24=27", 2)) { $count++ }
if ($d.Content.Find.Execute("15+38=53", $true, $false, $false, $false, $false, $true, 1, $false, "62-7=55", 2)) { $count++ }
if ($d.Content.Find.Execute("86-19=67", $true, $false, $false, $false, $false, $true, 1, $false, "17+59=76", 2)) { $count++ }
if ($d.Content.Find.Execute("8+58=66", $true, $false, $false, $false, $false, $true, 1, $false, "6+88=94", 2)) { $count++ }

Write-Output "Replaced $count / 101 items"
